# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Apio" (Feria Lagunitas de Puerto Montt)
# at the top of the data block (row 321), pushing the existing rows down by
# one (321-413 -> 322-414) and growing the sheet to A1:R414.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 321:413 down one row, creating a blank row 321.
$ws.Rows("321:321").Insert()

# Populate the new row 321 with the latest weekly record.
$ws.Range("A321").Value = 4
$ws.Range("B321").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C321").Value = "Los Lagos"
$ws.Range("D321").Value = 44988
$ws.Range("E321").Value = 10
$ws.Range("F321").Value = 100112017
$ws.Range("G321").Value = "Apio"
$ws.Range("H321").Value = "Americana (o)"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 50
$ws.Range("K321").Value = 12000
$ws.Range("L321").Value = 12000
$ws.Range("M321").Value = 12000
$ws.Range("N321").Value = "$/docena de matas"
$ws.Range("O321").Value = "Región de Coquimbo"
$ws.Range("P321").Value = 2000
$ws.Range("Q321").Value = 6
$ws.Range("R321").Value = "Hortaliza"
